$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '42.477.45'
Set-TextValue "E2" '  +0.98%  '
Set-TextValue "D3" '2.300.87'
Set-TextValue "E3" '  -0.28%  '
Set-TextValue "E4" '  +0.14%  '
Set-TextValue "D5" '316.72'
Set-TextValue "E5" '  +1.45%  '
Set-TextValue "D6" '104.29'
Set-TextValue "E6" '  -1.43%  '
Set-TextValue "E8" '  +0.20%  '
Set-TextValue "D9" '0.611'
Set-TextValue "E9" '  +0.33%  '
Set-TextValue "D10" '39.99'
Set-TextValue "E10" '  -0.80%  '
Set-TextValue "D11" '0.0909'
Set-TextValue "E11" '  -0.69%  '
Set-TextValue "D12" '8.35'
Set-TextValue "E12" '  +0.76%  '
Set-TextValue "D13" '0.106'
Set-TextValue "E13" '  +0.44%  '
Set-TextValue "D14" '0.965'
Set-TextValue "E14" '  -1.24%  '
Set-TextValue "D15" '15.34'
Set-TextValue "E15" '  -1.42%  '
Set-TextValue "D16" '2.652.31'
Set-TextValue "E16" '  -0.13%  '
Set-TextValue "D17" '2.307.54'
Set-TextValue "E17" '  -0.06%  '
Set-TextValue "D18" '42.462.09'
Set-TextValue "E18" '  +0.74%  '
Set-TextValue "D19" '7.46'
Set-TextValue "E19" '  -2.67%  '
Set-TextValue "E20" '  +0.86%  '
Set-TextValue "D21" '73.31'
Set-TextValue "E21" '  -1.72%  '
Set-TextValue "B22" 'BitcoinCash'
Set-TextValue "C22" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D22" '276.45'
Set-TextValue "E22" '  +6.38%  '
Set-TextValue "B23" 'PancakeSwap'
Set-TextValue "C23" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D23" '3.53'
Set-TextValue "E23" '  +1.63%  '
Set-TextValue "D24" '11.14'
Set-TextValue "E24" '  +19.46%  '
Set-TextValue "D25" '2.27'
Set-TextValue "E25" '  -1.07%  '
Set-TextValue "E26" '  -0.40%  '
Set-TextValue "E28" '  +3.31%  '
Set-TextValue "D29" '22.76'
Set-TextValue "E29" '  -0.02%  '
Set-TextValue "D30" '35.73'
Set-TextValue "E30" '  -0.35%  '
Set-TextValue "D31" '165.15'
Set-TextValue "E31" '  +0.77%  '
Set-TextValue "D32" '0.0871'
Set-TextValue "E32" '  -3.00%  '
Set-TextValue "D33" '5.88'
Set-TextValue "E33" '  +0.44%  '
Set-TextValue "E34" '  +3.93%  '
Set-TextValue "B35" 'Kaspa'
Set-TextValue "C35" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D35" '0.117'
Set-TextValue "E35" '  -0.81%  '
Set-TextValue "B36" 'WEMIXToken'
Set-TextValue "C36" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D36" '2.59'
Set-TextValue "E36" '  -11.55%  '
Set-TextValue "D37" '0.0369'
Set-TextValue "E37" '  +4.78%  '
Set-TextValue "D38" '4.58'
Set-TextValue "E38" '  +1.22%  '
Set-TextValue "E39" '  +3.69%  '
Set-TextValue "E40" '  -0.75%  '
Set-TextValue "D41" '1.49'
Set-TextValue "E41" '  +1.78%  '
Set-TextValue "D42" '69.70'
Set-TextValue "E42" '  -3.23%  '
Set-TextValue "D43" '0.227'
Set-TextValue "E43" '  -0.50%  '
Set-TextValue "D44" '94.36'
Set-TextValue "E44" '  -4.13%  '
Set-TextValue "E45" '  +0.11%  '
Set-TextValue "D46" '81.75'
Set-TextValue "E46" '  +9.91%  '
Set-TextValue "D47" '12.06'
Set-TextValue "E47" '  -1.94%  '
Set-TextValue "D48" '113.09'
Set-TextValue "E48" '  +0.46%  '
Set-TextValue "D49" '8.91'
Set-TextValue "E49" '  -1.35%  '
Set-TextValue "D50" '1.594.01'
Set-TextValue "E50" '  +2.25%  '
Set-TextValue "D51" '5.18'
Set-TextValue "E51" '  -2.93%  '
